# Update forest data - 2026-02-02 12:30
#
# "New" sheet currently holds 6 rows (rows 2-7) that are now stale; they
# get appended to the bottom of "Previously added" (rows 448-453, with
# their own hyperlinks), and "New" is refreshed with 7 freshly-scraped
# rows (and their hyperlinks).

$wb    = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item("Previously added")
$wsNew = $wb.Worksheets.Item("New")

# Stable reference cells (untouched by this script) carrying the exact
# per-column style ("link" / plain-text / date) that every data row on
# both sheets already uses, e.g. row 447 - the last pre-existing row of
# "Previously added".
$styleA = $wsOld.Range("A447")
$styleB = $wsOld.Range("B447")
$styleC = $wsOld.Range("C447")
$styleD = $wsOld.Range("D447")
$styleE = $wsOld.Range("E447")
$styleF = $wsOld.Range("F447")

# Helper: force a value to be written as TEXT (shared string), even if it
# looks numeric (e.g. a cadastre number). The leading "'" forces text
# typing; the caller re-stamps the correct style afterwards (Hyperlinks.Add
# and the text trick both leave transient/incorrect styles behind).
function Set-TextValue {
    param($cell, [string]$val)
    $cell.Value = "'" + $val
}

function Restore-Style {
    param($range, $styleSrc)
    $styleSrc.Copy()
    $range.PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Step 1: move the 6 existing "New" rows (2-7) down to the bottom of the
# "Previously added" sheet, as rows 448-453, carrying their hyperlinks.
# ---------------------------------------------------------------------

$firstOldRow = 448
$srcRange = $wsNew.Range("A2:F7")
$dstCell  = $wsOld.Range("A" + $firstOldRow)
$srcRange.Copy($dstCell)

$oldLinks = @(
    "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/stalbes-pag/hkpmm.html",
    "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/zosenu-pag/chebn.html",
    "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/bebrenes-pag/hljhk.html",
    "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kalniesu-pag/bebfkn.html",
    "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/indras-pag/beplxd.html",
    "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/silmalas-pag/ljeii.html"
)

for ($i = 0; $i -lt $oldLinks.Count; $i++) {
    $row = $firstOldRow + $i
    $cell = $wsOld.Range("A" + $row)
    $wsOld.Hyperlinks.Add($cell, $oldLinks[$i]) | Out-Null
}

$lastOldRow = $firstOldRow + $oldLinks.Count - 1
# Hyperlinks.Add silently reskins cells with the built-in "Hyperlink"
# style; restore the sheet's normal link-column style afterwards.
Restore-Style $wsOld.Range("A" + $firstOldRow + ":A" + $lastOldRow) $styleA

# ---------------------------------------------------------------------
# Step 2: wipe the old rows out of "New" and replace them with the 7
# freshly scraped rows.
# ---------------------------------------------------------------------

$wsNew.Range("A2:F7").ClearContents()
$wsNew.Hyperlinks.Delete()

# Columns: A=link, B=price, C=districtText, D=areaText, E=cadastreText, F=date
$newRows = @(
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/aizkraukle-and-reg/plavinas/mfgbj.html";
       B="6 000 €";   C="Aizkraukle un raj."; D="1 ha.";     E="32420090035"; F=46055.46388888889 },
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/jaunaluksnes-pag/inghf.html";
       B="33 000 €";  C="Alūksne un raj.";    D="2.30 ha.";  E="36560130028"; F=46055.36041666666 },
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/izvaltas-pag/hmlmn.html";
       B="120 000 €"; C="Krāslava un raj.";   D="19 ha.";    E="60640020120"; F=46054.705555555556 },
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/kurmales-pag/chgio.html";
       B="13 900 €";  C="Kuldīga un raj.";    D="2 ha.";     E="";            F=46054.82986111111 },
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/blontu-pag/dioce.html";
       B="41 000 €";  C="Ludza un raj.";      D="4 ha.";     E="68440050028"; F=46052.674305555556 },
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/ogre-and-reg/lielvardes-l-t/lhlxf.html";
       B="28 000 €";  C="Ogre un raj.";       D="8 ha.";     E="74840070028"; F=46053.77916666667 },
    @{ A="https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/saldus/kgkjn.html";
       B="6 000 €";   C="Saldus un raj.";     D="1.35 ha.";  E="84480060140"; F=46055.40763888889 }
)

$firstNewRow = 2
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $firstNewRow + $i
    $data = $newRows[$i]

    $cellA = $wsNew.Range("A" + $row)
    $cellB = $wsNew.Range("B" + $row)
    $cellC = $wsNew.Range("C" + $row)
    $cellD = $wsNew.Range("D" + $row)
    $cellE = $wsNew.Range("E" + $row)
    $cellF = $wsNew.Range("F" + $row)

    Set-TextValue $cellA $data.A
    Set-TextValue $cellB $data.B
    Set-TextValue $cellC $data.C
    Set-TextValue $cellD $data.D
    Set-TextValue $cellE $data.E
    $cellF.Value = $data.F

    $wsNew.Hyperlinks.Add($cellA, $data.A) | Out-Null
}

$lastNewRow = $firstNewRow + $newRows.Count - 1
Restore-Style $wsNew.Range("A" + $firstNewRow + ":A" + $lastNewRow) $styleA
Restore-Style $wsNew.Range("B" + $firstNewRow + ":B" + $lastNewRow) $styleB
Restore-Style $wsNew.Range("C" + $firstNewRow + ":C" + $lastNewRow) $styleC
Restore-Style $wsNew.Range("D" + $firstNewRow + ":D" + $lastNewRow) $styleD
Restore-Style $wsNew.Range("E" + $firstNewRow + ":E" + $lastNewRow) $styleE
Restore-Style $wsNew.Range("F" + $firstNewRow + ":F" + $lastNewRow) $styleF
